# Jerico Tallorin Atienzo Q0544 — "adding averages and more checks"
#
# 1. Training Dashboard: refresh the "PERIOD TO EXPIRE" counters and
#    "LAST UPDATE" dates (report re-run 8 days later).
# 2. Exam Dashboard: clearer remark text + a bit more breathing room in
#    the COMMENTS column.
# 3. Make the title banner and table header text bold/white on both
#    sheets (the title banner also drops back to the default font size).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Training Dashboard")
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------
# Training Dashboard — updated figures
# ---------------------------------------------------------------------
$ws1.Range("H3").Value = 647
$ws1.Range("I3").NumberFormat = "@"
$ws1.Range("I3").Value = "16-Sep-2025"
$ws1.Range("I3").NumberFormat = "general"

$ws1.Range("H4").Value = 423
$ws1.Range("I4").NumberFormat = "@"
$ws1.Range("I4").Value = "16-Sep-2025"
$ws1.Range("I4").NumberFormat = "general"

$ws1.Range("H5").Value = 182
$ws1.Range("I5").NumberFormat = "@"
$ws1.Range("I5").Value = "16-Sep-2025"
$ws1.Range("I5").NumberFormat = "general"

# ---------------------------------------------------------------------
# Exam Dashboard — remark wording + wider COMMENTS column
# ---------------------------------------------------------------------
$ws2.Range("E3").Value = "date is valid"
$ws2.Range("E4").Value = "date is valid"
$ws2.Columns.Item(5).ColumnWidth = 14.166666666666666

# ---------------------------------------------------------------------
# Title banner (A1) + header row (row 2) formatting, both sheets
# ---------------------------------------------------------------------
$ws1.Range("A1").Font.Size = 11
$ws1.Range("A1").Font.Color = 16777215
$ws1.Range("A2:K2").Font.Size = 11
$ws1.Range("A2:K2").Font.Color = 16777215
$ws1.Range("A2:K2").Font.Bold = $true

$ws2.Range("A1").Font.Size = 11
$ws2.Range("A1").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Size = 11
$ws2.Range("A2:G2").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Bold = $true
